$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A55").Value = 'plate; dish'
$ws.Range("B55").Value = '皿|さら'
$ws.Range("A56").Value = 'paper plate'
$ws.Range("B56").Value = '紙皿|かみざら'
$ws.Range("A57").Value = 'ashtray'
$ws.Range("B57").Value = '灰皿|はいざら'
$ws.Range("A58").Value = 'voice'
$ws.Range("B58").Value = '声|こえ'
$ws.Range("A59").Value = 'phonetics'
$ws.Range("B59").Value = '音声学|おんせいがく'
$ws.Range("A60").Value = 'onomatopoeia'
$ws.Range("B60").Value = '擬声語|ぎせいご'
$ws.Range("A61").Value = 'voice actor'
$ws.Range("B61").Value = '声優|せいゆう'
$ws.Range("A62").Value = 'Japanese tea'
$ws.Range("B62").Value = 'お茶|おちゃ'
$ws.Range("A63").Value = 'teahouse'
$ws.Range("B63").Value = '茶店|ちゃみせ'
$ws.Range("A64").Value = 'black tea'
$ws.Range("B64").Value = '紅茶|こうちゃ'
$ws.Range("A65").Value = 'brown'
$ws.Range("B65").Value = '茶色|ちゃいろ'
$ws.Range("A66").Value = 'cafe'
$ws.Range("B66").Value = '喫茶店|きっさてん'
$ws.Range("A67").Value = '(something) stops'
$ws.Range("B67").Value = '止まる|とまる'
$ws.Range("A68").Value = 'to cancel'
$ws.Range("B68").Value = '中止する|ちゅうしする'
$ws.Range("A69").Value = 'to prohibit'
$ws.Range("B69").Value = '禁止する|きんしする'
$ws.Range("A70").Value = 'one sheet'
$ws.Range("B70").Value = '一枚|いちまい'
$ws.Range("A71").Value = 'number of flat things'
$ws.Range("B71").Value = '枚数|まいすう'
$ws.Range("A72").Value = 'Three ryoo'
$ws.Range("B72").Value = '三両|さんりょう'
$ws.Range("A73").Value = 'parents'
$ws.Range("B73").Value = '両親|りょうしん'
$ws.Range("A74").Value = 'both hands'
$ws.Range("B74").Value = '両手|りょうて'
$ws.Range("A75").Value = 'exchange'
$ws.Range("B75").Value = '両替|りょうがえ'
$ws.Range("A76").Value = 'both sides'
$ws.Range("B76").Value = '両方|りょうほう'
$ws.Range("A77").Value = 'impossible'
$ws.Range("B77").Value = '無理な|むりな'
$ws.Range("A78").Value = 'wasteful'
$ws.Range("B78").Value = '無駄な|むだな'
$ws.Range("A79").Value = 'free of charge'
$ws.Range("B79").Value = '無料|むりょう'
$ws.Range("A80").Value = 'rude'
$ws.Range("B80").Value = '無礼な|ぶれいな'
$ws.Range("A81").Value = 'there is no...'
$ws.Range("B81").Value = '無い|ない'
$ws.Range("A82").Value = 'to pay'
$ws.Range("B82").Value = '払う|はらう'
$ws.Range("A83").Value = 'payment'
$ws.Range("B83").Value = '支払い|しはらい'
$ws.Range("A84").Value = 'refund'
$ws.Range("B84").Value = '払い戻し|はらいもどし'
$ws.Range("A85").Value = 'payment in installments'
$ws.Range("B85").Value = '分割払い|ぶんかつばらい'
$ws.Range("A86").Value = 'heart; mind'
$ws.Range("B86").Value = '心|こころ'
$ws.Range("A87").Value = 'to worry'
$ws.Range("B87").Value = '心配する|しんぱいする'
$ws.Range("A88").Value = 'enthusiastic'
$ws.Range("B88").Value = '熱心な|ねっしんな'
$ws.Range("A89").Value = 'safe'
$ws.Range("B89").Value = '安心な|あんしんな'
$ws.Range("A90").Value = 'curiosity'
$ws.Range("B90").Value = '好奇心|こうきしん'
$ws.Range("A91").Value = 'to laugh'
$ws.Range("B91").Value = '笑う|わらう'
$ws.Range("A92").Value = 'smile; smiling face'
$ws.Range("B92").Value = '笑顔|えがお'
$ws.Range("A93").Value = 'to smile'
$ws.Range("B93").Value = '微笑む|ほほえむ'
$ws.Range("A94").Value = 'to burst into laughter'
$ws.Range("B94").Value = '爆笑する|ばくしょうする'
$ws.Range("A95").Value = 'definitely'
$ws.Range("B95").Value = '絶対に|ぜったいに'
$ws.Range("A96").Value = 'to die out'
$ws.Range("B96").Value = '絶える|たえる'
$ws.Range("A97").Value = 'extinction'
$ws.Range("B97").Value = '絶滅|ぜつめつ'
$ws.Range("A98").Value = 'to faint'
$ws.Range("B98").Value = '気絶する|きぜつする'
$ws.Range("A99").Value = 'despair'
$ws.Range("B99").Value = '絶望|ぜつぼう'
$ws.Range("A100").Value = 'to oppose'
$ws.Range("B100").Value = '反対する|はんたいする'
$ws.Range("A101").Value = 'Japan versus China'
$ws.Range("B101").Value = '日本対中国|にほんたいちゅうごく'
$ws.Range("A102").Value = 'a pair'
$ws.Range("B102").Value = '一対|いっつい'
$ws.Range("A103").Value = 'painful'
$ws.Range("B103").Value = '痛い|いたい'
$ws.Range("A104").Value = 'painkiller'
$ws.Range("B104").Value = '痛み止め|いたみどめ'
$ws.Range("A105").Value = 'headache'
$ws.Range("B105").Value = '頭痛|ずつう'
$ws.Range("A106").Value = 'stomachache'
$ws.Range("B106").Value = '腹痛|ふくつう'
$ws.Range("A107").Value = 'lower back pain'
$ws.Range("B107").Value = '腰痛|ようつう'
$ws.Range("A108").Value = 'the worst'
$ws.Range("B108").Value = '最悪|さいあく'
$ws.Range("A109").Value = 'recently'
$ws.Range("B109").Value = '最近|さいきん'
$ws.Range("A110").Value = 'the best'
$ws.Range("B110").Value = '最高|さいこう'
$ws.Range("A111").Value = 'the latest'
$ws.Range("B111").Value = '最新|さいしん'
$ws.Range("A112").Value = 'lastly'
$ws.Range("B112").Value = '最後に|さいごに'
$ws.Range("A113").Value = 'most'
$ws.Range("B113").Value = '最も|もっとも'
$ws.Range("A114").Value = 'to continue'
$ws.Range("B114").Value = '続ける|つづける'
$ws.Range("A115").Value = 'procedure'
$ws.Range("B115").Value = '手続き|てつづき'
$ws.Range("A116").Value = 'to inherit'
$ws.Range("B116").Value = '相続する|そうぞくする'
$ws.Range("A117").Value = 'serial TV drama'
$ws.Range("B117").Value = '連続ドラマ|れんぞくドラマ'
